# Updates cryptocurrency price/volume data per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "93.517.93"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.441.37"
$ws.Range("E3").Value = "  +4.44%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +10.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.392"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.00%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +12.69%  "
$ws.Range("D11").Value = "3.439.24"
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.199"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.08%  "
$ws.Range("D15").Value = "4.099.95"
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("D16").Value = "93.383.20"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000249"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.74%  "
$ws.Range("D19").Value = "3.445.12"
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.491"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.67%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "504.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.97%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.01%  "
$ws.Range("E26").Value = "  +3.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.72%  "
$ws.Range("D29").Value = "3.624.88"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.03%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.50%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.139"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.67%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.179"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.554"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "570.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +3.69%  "
$ws.Range("E43").Value = "  +7.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("E45").Value = "  +11.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.02%  "
$ws.Range("E47").Value = "  +5.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +7.41%  "
